# Natmi following Dr Hou advice
# Update the LR-pairs_lrc2p Fn1-Plaur sheet: Ligand-expressing cells (E) and
# Receptor-expressing cells (K) change from 1 to 3, and all downstream
# expression/specificity metrics (G,H,I,J,M,N,O,P,Q,R,S,T) are recomputed
# accordingly for data rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 56.98117766666667
$ws.Range("N2").Value = 170.943533
$ws.Range("O2").Value = 0.952030123851636
$ws.Range("P2").Value = 0.9520301238516359
$ws.Range("Q2").Value = 1540.708687704573
$ws.Range("R2").Value = 13866.37818934116
$ws.Range("S2").Value = 0.06755784949041201
$ws.Range("T2").Value = 0.06755784949041201

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.516719
$ws.Range("N3").Value = 7.550157
$ws.Range("O3").Value = 0.04204883786863874
$ws.Range("P3").Value = 0.04204883786863874
$ws.Range("Q3").Value = 68.04932763050766
$ws.Range("R3").Value = 612.443948674569
$ws.Range("S3").Value = 0.002983864679074935
$ws.Range("T3").Value = 0.002983864679074936

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("N4").Value = 1.063163
$ws.Range("O4").Value = 0.005921038279725251
$ws.Range("P4").Value = 0.005921038279725251
$ws.Range("Q4").Value = 9.582254688430108
$ws.Range("R4").Value = 86.24029219587098
$ws.Range("S4").Value = 0.0004201680208503406
$ws.Range("T4").Value = 0.0004201680208503406

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 56.98117766666667
$ws.Range("N5").Value = 170.943533
$ws.Range("O5").Value = 0.952030123851636
$ws.Range("P5").Value = 0.9520301238516359
$ws.Range("Q5").Value = 19690.7906336612
$ws.Range("R5").Value = 177217.1157029508
$ws.Range("S5").Value = 0.8634127142866952
$ws.Range("T5").Value = 0.8634127142866952

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.516719
$ws.Range("N6").Value = 7.550157
$ws.Range("O6").Value = 0.04204883786863874
$ws.Range("P6").Value = 0.04204883786863874
$ws.Range("Q6").Value = 869.693975134301
$ws.Range("R6").Value = 7827.245776208709
$ws.Range("S6").Value = 0.03813482402203944
$ws.Range("T6").Value = 0.03813482402203944

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("N7").Value = 1.063163
$ws.Range("O7").Value = 0.005921038279725251
$ws.Range("P7").Value = 0.005921038279725251
$ws.Range("Q7").Value = 122.4645336097923
$ws.Range("R7").Value = 1102.180802488131
$ws.Range("S7").Value = 0.005369892826300633
$ws.Range("T7").Value = 0.005369892826300634

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.98117766666667
$ws.Range("N8").Value = 170.943533
$ws.Range("O8").Value = 0.952030123851636
$ws.Range("P8").Value = 0.9520301238516359
$ws.Range("Q8").Value = 480.2794554712365
$ws.Range("R8").Value = 4322.515099241128
$ws.Range("S8").Value = 0.02105956007452876
$ws.Range("T8").Value = 0.02105956007452876

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.516719
$ws.Range("N9").Value = 7.550157
$ws.Range("O9").Value = 0.04204883786863874
$ws.Range("P9").Value = 0.04204883786863874
$ws.Range("Q9").Value = 21.21276674843467
$ws.Range("R9").Value = 190.914900735912
$ws.Range("S9").Value = 0.0009301491675243651
$ws.Range("T9").Value = 0.0009301491675243653

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("N10").Value = 1.063163
$ws.Range("O10").Value = 0.005921038279725251
$ws.Range("P10").Value = 0.005921038279725251
$ws.Range("Q10").Value = 2.987041029023111
$ws.Range("R10").Value = 26.883369261208
$ws.Range("S10").Value = 0.0001309774325742771
$ws.Range("T10").Value = 0.0001309774325742771
